$d = $word.ActiveDocument

# --- 1. Add the new YouTube hyperlink to the (currently empty) third list item ---
$linkPara = $d.Paragraphs(4)
$link = $d.Hyperlinks.Add($linkPara.Range, "https://www.youtube.com/watch?v=RWXKysImabs", `
    [Type]::Missing, [Type]::Missing, "https://www.youtube.com/watch?v=RWXKysImabs")

# --- 2. Turn the "Single Page Architecture:" paragraph into a new (empty) list item ---
$spaPara = $d.Paragraphs(5)

# 2a. Insert a brand-new blank paragraph right after it *before* changing its
#     formatting, so the new paragraph keeps the original plain/Normal formatting.
$spaPara.Range.InsertParagraphAfter()

# 2b. Remove the paragraph's text (but keep the paragraph mark) ...
$textRange = $d.Range($spaPara.Range.Start, $spaPara.Range.End - 1)
$textRange.Delete()

# 2c. ... and give it the same list formatting (List Paragraph style, same
#      numbered list) as the other bullet items above it.
$spaPara.Style = "List Paragraph"
$listTemplate = $linkPara.Range.ListFormat.ListTemplate
$spaPara.Range.ListFormat.ApplyListTemplateWithLevel($listTemplate, $true)

Write-Output "done"
